$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 406, pushing the existing row 406 (and
# everything below it) down by one. This mirrors the diff, which shows
# every record from the old row 406 through row 500 reappearing one row
# lower (407-501), with a brand-new weekly price record occupying the
# freshly inserted row 406.
$ws.Rows("406").Insert()

# Populate the new row 406 with the new weekly record. Columns A, B, C,
# E, F, G, N, Q, R hold the same constant values as every other data row
# in this sheet.
$ws.Range("A406").Value = 4
$ws.Range("B406").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C406").Value = "Los Lagos"
$ws.Range("D406").Value = 44785
$ws.Range("E406").Value = 10
$ws.Range("F406").Value = 100112006
$ws.Range("G406").Value = "Repollo"
$ws.Range("H406").Value = "Crespo record"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 1400
$ws.Range("K406").Value = 2300
$ws.Range("L406").Value = 2400
$ws.Range("M406").Value = 2350
$ws.Range("N406").Value = "$/unidad"
$ws.Range("O406").Value = "Región Metropolitana"
$ws.Range("P406").Value = 2350
$ws.Range("Q406").Value = 1
$ws.Range("R406").Value = "Hortaliza"
